$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.246.60'
$ws.Range("E2").Value = '  +3.72%  '
$ws.Range("D3").Value = '2.062.03'
$ws.Range("E3").Value = '  +3.32%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("E6").Value = '  +2.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.72'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.23%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0810'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.28%  '
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '2.364.28'
$ws.Range("E12").Value = '  +3.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.00%  '
$ws.Range("E14").Value = '  +3.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.755'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.14%  '
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("D17").Value = '2.053.72'
$ws.Range("E17").Value = '  +3.01%  '
$ws.Range("D18").Value = '38.123.34'
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("E25").Value = '  +5.23%  '
$ws.Range("E26").Value = '  +3.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("E28").Value = '  +9.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("E33").Value = '  +6.28%  '
$ws.Range("E34").Value = '  +1.77%  '
$ws.Range("E35").Value = '  +7.54%  '
$ws.Range("E36").Value = '  +1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +15.74%  '
$ws.Range("E38").Value = '  +7.03%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.89%  '
$ws.Range("D42").Value = '1.482.92'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0949'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.73%  '
$ws.Range("E45").Value = '  +4.00%  '
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("E47").Value = '  +18.03%  '
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = '2.251.22'
$ws.Range("E51").Value = '  +3.37%  '
